# Append new order rows (112-121) to the "Orders" sheet, and extend the
# tracking-number string on the "Summary" sheet's G2 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New rows to append after the existing last row (111).
# Columns used: A = PackageID, C = FlowerName, F = Number.
# All of these are stored as text (same convention as every other row in
# this sheet - numeric-looking values are kept as text, not numbers).
$newRows = @(
    @{ Row = 112; A = $null; C = "532_灯苔_undefined_undefined_1bunch"; F = "5" },
    @{ Row = 113; A = "2";   C = "144_高原红_High Plateau Red_Rosa rugosa Thunb._20stems"; F = "12" },
    @{ Row = 114; A = $null; C = "268_猩红泡泡_spray red_Rosa rugosa Thunb._10stems"; F = "11" },
    @{ Row = 115; A = $null; C = "259_诺拉_Nora_Rosa rugosa Thunb._10stems"; F = "4" },
    @{ Row = 116; A = $null; C = "149_骄傲_Proud_Rosa rugosa Thunb._20stems"; F = "14" },
    @{ Row = 117; A = $null; C = "624_多丁白_undefined_undefined_1bunch"; F = "10" },
    @{ Row = 118; A = $null; C = "608_康乃馨笑颜_undefined_undefined_20stems"; F = "15" },
    @{ Row = 119; A = "3";   C = "105_绣球莫奈蓝_Hydrangea Monet Blue_Hydrangea L._1stem"; F = "20" },
    @{ Row = 120; A = $null; C = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"; F = "20" },
    @{ Row = 121; A = $null; C = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"; F = $null }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    if ($null -ne $r.A) {
        $cell = $ws.Cells.Item($rowNum, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $r.A
    }

    if ($null -ne $r.C) {
        $cell = $ws.Cells.Item($rowNum, 3)
        $cell.NumberFormat = "@"
        $cell.Value = $r.C
    }

    if ($null -ne $r.F) {
        $cell = $ws.Cells.Item($rowNum, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $r.F
    }
}

# Extend the concatenated tracking-number string on Summary!G2.
# Build the new value via [string]::Concat (not "+", which this host's
# PowerShell coerces to numeric addition for two numeric-looking strings)
# and keep the cell text-formatted so the huge numeral isn't rounded into
# scientific notation when written back.
$summary = $wb.Worksheets.Item("Summary")
$g2 = $summary.Range("G2")
$g2NewValue = [string]::Concat($g2.Value(), "51211414101520200")
$g2.NumberFormat = "@"
$g2.Value = $g2NewValue
